$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Header row (row 1) gains a new column H (array_size) ---
$ws.Range("A1").Value = "type"
$ws.Range("B1").Value = "name"
$ws.Range("C1").Value = "topic"
$ws.Range("D1").Value = "source"
$ws.Range("E1").Value = "module"
$ws.Range("F1").Value = "dtype"
$ws.Range("G1").Value = "value_units"

# New H1 header cell - carry over the header formatting from G1
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)
$ws.Range("H1").Value = "array_size"

# --- Row 2: NXdetector / detector, clear the topic/source/module/dtype/value_units columns ---
$ws.Range("A2").Value = "NXdetector"
$ws.Range("B2").Value = "detector"
$ws.Range("C2:G2").ClearContents()

# --- Row 3: NXsample / sample, clear the rest ---
$ws.Range("A3").Value = "NXsample"
$ws.Range("B3").Value = "sample"
$ws.Range("C3:G3").ClearContents()

# --- Row 4: NXmonitor / control, clear the rest ---
$ws.Range("A4").Value = "NXmonitor"
$ws.Range("B4").Value = "control"
$ws.Range("C4:G4").ClearContents()

# --- Row 5: NXdata / data / odin_topic / image_source / ADAr / uint32 / (blank) / array_size value ---
$ws.Range("A5").Value = "NXdata"
$ws.Range("B5").Value = "data"
$ws.Range("C5").Value = "odin_topic"
$ws.Range("D5").Value = "image_source"
$ws.Range("E5").Value = "ADAr"
$ws.Range("F5").Value = "uint32"

# New H5 value cell - carry over the data-row formatting from G5
$ws.Range("G5").Copy()
$ws.Range("H5").PasteSpecial(-4122)
$ws.Range("H5").Value = "480, 290, 3"

# --- View tweaks: scroll so column B is the left edge, select H7 ---
$excel.ActiveWindow.ScrollColumn = 2
$ws.Range("H7").Select() | Out-Null
